$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that are no longer present in the updated tax summary:
#  - rows that collapsed to ~0 / were dropped from the table entirely
#  - rows 23, 26 & 27 (the "cannot be assigned", "unclassified" and
#    "Viruses" rows) which get folded into a single combined
#    "Unclassified" row appended at the end.
# Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(27,26,25,23,22,20,18,17,8,5,2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Append the new combined "Unclassified" row (replaces the three removed
# unclassified-like rows with summed Summer/Winter/Spring percentages).
$ws.Range("A17").Value = "Unclassified"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("B17").Value = "NA"
$ws.Range("C17").Value = "NA"
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "50.742 ± 2.174"
$ws.Range("F17").Value = "26.972 ± 6.522"
$ws.Range("G17").Value = "54.182 ± 1.996"

# Match the saved selection state.
$ws.Range("G17").Select()
